# Refresh the cryptos list: update each coin's Price (column D) and
# Volume(1h) change (column E) to the latest scraped values.
#
# Values are written with a leading apostrophe (quote-prefix) so Excel
# stores them as literal text instead of auto-coercing numeric-looking
# strings (e.g. "225.61", "4.610", "1.726.96") into floating point
# numbers -- which would silently drop significant trailing zeros and
# lose the thousands-grouped "27.533.67" style formatting used
# throughout this sheet. The cell Style is then reset to "Normal" so no
# stray number-format/quote-prefix style index stays attached to the
# cell (matches the source workbook, where these cells carry no `s`
# attribute).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.533.67'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +5.45%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''1.726.96'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +4.55%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.15%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''225.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +3.29%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''0.5361'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +2.86%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D8").Value = '''0.2667'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +0.87%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.06602'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +4.16%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''21.74'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +6.71%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.07723'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +0.40%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''4.610'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -0.26%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''1.732.29'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +3.88%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''1.963.68'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +4.47%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = '''  +4.43%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''0.0₅8293'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +1.70%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''67.98'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +4.00%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''27.558.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = '''220.16'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +15.44%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''1.004'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +0.09%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''4.735'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +2.33%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''10.65'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +1.72%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''6.099'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +2.79%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = '''  +0.06%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''  +2.53%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''1.710'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +13.11%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''0.1234'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +3.72%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''7.416'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +2.75%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '''  +4.78%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''0.05570'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +1.47%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = '''  +2.68%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D33").Value = '''3.459'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +2.94%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''1.660'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +6.50%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''0.9622'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +1.39%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  +1.47%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''2.430'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +1.31%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''0.5949'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +5.49%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''0.01652'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +4.74%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = '''  +1.21%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '''  +3.06%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''1.057.73'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +2.81%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  +0.08%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''101.42'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +0.41%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''1.870.66'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = '''0.0₈115'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +5.53%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''59.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +2.41%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''8.202'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +2.92%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  +2.30%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''  +0.15%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''0.06545'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +12.54%  '
$ws.Range("E51").Style = "Normal"
